$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Fix Fusion's quantity from 3 to 8
$table.Cell(2, 3).Range.Text = "8"

# New vehicle rows to append at the bottom of the table
$newRows = @(
    @("vehicle", "Dondai",  "4", "30000.0", "2009", "4", "5"),
    @("vehicle", "Civic",   "6", "24000.0", "2011", "4", "8"),
    @("vehicle", "Charger", "8", "20004.0", "2004", "2", "20"),
    @("vehicle", "Fiesta",  "2", "14000.0", "2006", "3", "23")
)

foreach ($rowData in $newRows) {
    $table.Rows.Add() | Out-Null
    $r = $table.Rows.Count
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $table.Cell($r, $c).Range.Text = $rowData[$c - 1]
    }
}
